$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5810.8887
$ws.Range("I76").Value = 5833.3335
$ws.Range("J76").Value = 5799.6665
$ws.Range("K76").Value = 5833.3335
$ws.Range("L76").Value = 5799.6665
$ws.Range("M76").Value = -5518.3335
$ws.Range("N76").Value = -6429.6665

$ws.Range("H79").Value = 5810.8887
$ws.Range("I79").Value = 5833.3335
$ws.Range("J79").Value = 5799.6665
$ws.Range("K79").Value = 5833.3335
$ws.Range("L79").Value = 5799.6665
$ws.Range("M79").Value = -4741.3335
$ws.Range("N79").Value = -7983.6665

$ws.Range("H98").Value = 5736
$ws.Range("I98").Value = 6249
$ws.Range("J98").Value = 606
$ws.Range("K98").Value = 6249
$ws.Range("L98").Value = 606
$ws.Range("M98").Value = -4751
$ws.Range("N98").Value = -3602

$ws.Range("H122").Value = 5736
$ws.Range("I122").Value = 6249
$ws.Range("J122").Value = 606
$ws.Range("K122").Value = 18747
$ws.Range("L122").Value = 1818
$ws.Range("M122").Value = -16297
$ws.Range("N122").Value = -6718

$ws.Range("H135").Value = 25641640
$ws.Range("I135").Value = 465.68967
$ws.Range("J135").Value = 100001050
$ws.Range("K135").Value = 4191.20703
$ws.Range("L135").Value = 900009450
$ws.Range("M135").Value = -1656.20703
$ws.Range("N135").Value = -900014520

$ws.Range("H137").Value = 1852.341
$ws.Range("I137").Value = 1823.75
$ws.Range("K137").Value = 5471.25
$ws.Range("M137").Value = -2921.25

$ws.Range("H138").Value = 2211.7087
$ws.Range("I138").Value = 1990.75
$ws.Range("J138").Value = 2267.8254
$ws.Range("K138").Value = 5972.25
$ws.Range("L138").Value = 6803.476200000001
$ws.Range("M138").Value = -832.25
$ws.Range("N138").Value = -17083.4762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 864.6429000000001
$ws.Range("I2").Value = 628.9545000000001
$ws.Range("J2").Value = 1728.8334
$ws.Range("K2").Value = 628.9545000000001
$ws.Range("L2").Value = 1728.8334
$ws.Range("M2").Value = -515.9545000000001
$ws.Range("N2").Value = -1954.8334

$ws.Range("H45").Value = 1255.4445
$ws.Range("I45").Value = 1203.1666
$ws.Range("J45").Value = 1360
$ws.Range("K45").Value = 1203.1666
$ws.Range("L45").Value = 1360
$ws.Range("M45").Value = -826.1666
$ws.Range("N45").Value = -2114

$ws.Range("H61").Value = 32259290
$ws.Range("I61").Value = 38462510
$ws.Range("K61").Value = 38462510
$ws.Range("M61").Value = -38462298

$ws.Range("H74").Value = 2917.6
$ws.Range("I74").Value = 2074.2
$ws.Range("K74").Value = 2074.2
$ws.Range("M74").Value = -1200.2

$ws.Range("H77").Value = 2917.6
$ws.Range("I77").Value = 2074.2
$ws.Range("K77").Value = 10371
$ws.Range("M77").Value = -6003

$ws.Range("H97").Value = 561.8461
$ws.Range("I97").Value = 546.1667
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 546.1667
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -50.16669999999999
$ws.Range("N97").Value = -1742

$ws.Range("H116").Value = 864.6429000000001
$ws.Range("I116").Value = 628.9545000000001
$ws.Range("J116").Value = 1728.8334
$ws.Range("K116").Value = 628.9545000000001
$ws.Range("L116").Value = 1728.8334
$ws.Range("M116").Value = 1665.0455
$ws.Range("N116").Value = -6316.8334

$ws.Range("H122").Value = 1550.75
$ws.Range("I122").Value = 1468.375
$ws.Range("K122").Value = 4405.125
$ws.Range("M122").Value = -1955.125

$ws.Range("H136").Value = 32259290
$ws.Range("I136").Value = 38462510
$ws.Range("K136").Value = 115387530
$ws.Range("M136").Value = -115384980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 864.6429000000001
$ws.Range("I3").Value = 628.9545000000001
$ws.Range("J3").Value = 1728.8334
$ws.Range("K3").Value = 628.9545000000001
$ws.Range("L3").Value = 1728.8334
$ws.Range("M3").Value = -514.9545000000001
$ws.Range("N3").Value = -1956.8334

$ws.Range("H105").Value = 1527.6923
$ws.Range("I105").Value = 1396.3636
$ws.Range("K105").Value = 1396.3636
$ws.Range("M105").Value = 350.6364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1065.4584
$ws.Range("I31").Value = 929.8889
$ws.Range("K31").Value = 929.8889
$ws.Range("M31").Value = -634.8889

$ws.Range("H34").Value = 1065.4584
$ws.Range("I34").Value = 929.8889
$ws.Range("K34").Value = 929.8889
$ws.Range("M34").Value = -727.8889

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 840.1053000000001
$ws.Range("I122").Value = 851.2308
$ws.Range("J122").Value = 816
$ws.Range("K122").Value = 2553.6924
$ws.Range("L122").Value = 2448
$ws.Range("M122").Value = -103.6923999999999
$ws.Range("N122").Value = -7348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 30349232
$ws.Range("J131").Value = 75936.39999999999
$ws.Range("L131").Value = 227809.2
$ws.Range("N131").Value = -237889.2

$ws.Range("H132").Value = 1205.4375
$ws.Range("I132").Value = 976.44446
$ws.Range("J132").Value = 1499.8572
$ws.Range("K132").Value = 8788.00014
$ws.Range("L132").Value = 13498.7148
$ws.Range("M132").Value = -6258.00014
$ws.Range("N132").Value = -18558.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7500000
$ws.Range("I10").Value = 7500000
$ws.Range("K10").Value = 7500000
$ws.Range("M10").Value = -7499831

$ws.Range("H70").Value = 26474342
$ws.Range("I70").Value = 17861072
$ws.Range("J70").Value = 66669600
$ws.Range("K70").Value = 17861072
$ws.Range("L70").Value = 66669600
$ws.Range("M70").Value = -17860802
$ws.Range("N70").Value = -66670140

$ws.Range("H73").Value = 26474342
$ws.Range("I73").Value = 17861072
$ws.Range("J73").Value = 66669600
$ws.Range("K73").Value = 17861072
$ws.Range("L73").Value = 66669600
$ws.Range("M73").Value = -17860136
$ws.Range("N73").Value = -66671472

$ws.Range("H122").Value = 5997
$ws.Range("I122").Value = 5997
$ws.Range("K122").Value = 17991
$ws.Range("M122").Value = -15541

$ws.Range("H132").Value = 3026.5667
$ws.Range("I132").Value = 2666.0417
$ws.Range("K132").Value = 7998.125100000001
$ws.Range("M132").Value = -5468.125100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 764.5454999999999
$ws.Range("I22").Value = 676.25
$ws.Range("K22").Value = 676.25
$ws.Range("M22").Value = -381.25

$ws.Range("H27").Value = 764.5454999999999
$ws.Range("I27").Value = 676.25
$ws.Range("K27").Value = 676.25
$ws.Range("M27").Value = -569.25

$ws.Range("H40").Value = 3472.6667
$ws.Range("I40").Value = 2410.111
$ws.Range("J40").Value = 5066.5
$ws.Range("K40").Value = 2410.111
$ws.Range("L40").Value = 5066.5
$ws.Range("M40").Value = -2274.111
$ws.Range("N40").Value = -5338.5

$ws.Range("H46").Value = 4133.3335
$ws.Range("J46").Value = 4615.385
$ws.Range("L46").Value = 4615.385
$ws.Range("N46").Value = -4991.385

$ws.Range("H82").Value = 1897.4
$ws.Range("I82").Value = 1861
$ws.Range("J82").Value = 2043
$ws.Range("K82").Value = 1861
$ws.Range("L82").Value = 2043
$ws.Range("M82").Value = -1500
$ws.Range("N82").Value = -2765

$ws.Range("H85").Value = 1897.4
$ws.Range("I85").Value = 1861
$ws.Range("J85").Value = 2043
$ws.Range("K85").Value = 1861
$ws.Range("L85").Value = 2043
$ws.Range("M85").Value = -613
$ws.Range("N85").Value = -4539

$ws.Range("H100").Value = 1575.75
$ws.Range("I100").Value = 1201.5
$ws.Range("K100").Value = 1201.5
$ws.Range("M100").Value = -660.5

$ws.Range("H122").Value = 41668252
$ws.Range("I122").Value = 62501376
$ws.Range("K122").Value = 187504128
$ws.Range("M122").Value = -187501678

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H20").Value = 86677.664
$ws.Range("J20").Value = 86677.664
$ws.Range("L20").Value = 86677.664
$ws.Range("N20").Value = -87157.664

$ws.Range("H107").Value = 571.6875
$ws.Range("J107").Value = 701.95
$ws.Range("L107").Value = 2105.85
$ws.Range("N107").Value = -5945.85

$ws.Range("H122").Value = 58825780
$ws.Range("I122").Value = 62502268
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 187506804
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -187504354
$ws.Range("N122").Value = -10870
